# "Add files via upload" / "Celulas coloridas agora!"
#
# Renames the two data sheets, widens the columns on the main sheet,
# colors the data rows of the main sheet according to the owning
# department (blue / orange / green), and populates the previously
# empty "Plan2" sheet with a bold+underlined header and the list of
# sectors that were removed from PADS_SECGS_Horas.xlsx.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename worksheets
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Siglas Completas(Atngias+Novas)"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "PADS_SECGS"

# ---------------------------------------------------------------
# 2. Sheet2 (PADS_SECGS): add header first (bold + underlined)
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "SETORES RETIRADOS DA PLANILHA: PADS_SECGS_Horas.xlsx"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Underline = $true

# ---------------------------------------------------------------
# 3. Sheet1: widen columns A and B (best-fit sized wider to fit
#    content) and color the rows by department
# ---------------------------------------------------------------

# ColumnWidth values are compensated so the persisted column width
# (which is ~0.8333 characters wider than the value assigned through
# this interop layer) lands on the intended widths.
$ws1.Columns.Item(1).ColumnWidth = 27.022135416666668
$ws1.Columns.Item(2).ColumnWidth = 59.022135416666664

# Rows 14-23 (Coordenadoria de Seguranca, Transporte e Apoio Administrativo) -> green
$ws1.Range("A14:B23").Interior.Color = 0xBDE4D7

# Rows 4-13 (Coordenadoria de InfraEstrutura Predial) -> orange
$ws1.Range("A4:B13").Interior.Color = 0xB5D5FC

# Rows 2-3 (SECGS, GABGS / Secretaria de Gestao de Servicos) -> blue
$ws1.Range("A2:B3").Interior.Color = 0xE5CDB9

# ---------------------------------------------------------------
# 4. Sheet2 (PADS_SECGS): list of removed sectors
# ---------------------------------------------------------------
$ws2.Range("A2").Value = "SECGS"
$ws2.Range("A3").Value = "GABGS"
$ws2.Range("A4").Value = "CIP"
$ws2.Range("A5").Value = "SAPRE"
$ws2.Range("A6").Value = "SMIC"
$ws2.Range("A7").Value = "SMIN"
$ws2.Range("A8").Value = "SOP"
$ws2.Range("A9").Value = "CSTA"
$ws2.Range("A10").Value = "SEXP"
$ws2.Range("A11").Value = "ST"
$ws2.Range("A12").Value = "SESEG"

$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 5. View state: selection on each sheet and zoom on sheet1
# ---------------------------------------------------------------
$ws2.Range("A2:A12").Select() | Out-Null

$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 115
$ws1.Range("G3").Select() | Out-Null
